$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shift the old row 24 ("30 channels duration") down to row 27, making room
# for the two new rows (RX_OFF / RX_ON duration) plus the new row 22 gap and
# the new ON/OFF ratio row.
$ws.Rows.Item(24).Insert()
$ws.Rows.Item(24).Insert()
$ws.Rows.Item(24).Insert()

# --- Row 25: RX_OFF duration ---
$ws.Range("A25").Value = "RX_OFF duration"
$ws.Range("A25").Style = "Good"
$ws.Range("B25").Value = 45
$ws.Range("B25").Style = "Good"
$ws.Range("C25").Value = "ms"

# --- Row 26: RX_ON duration ---
$ws.Range("A26").Value = "RX_ON duration"
$ws.Range("A26").Style = "Good"
$ws.Range("B26").Value = 18
$ws.Range("B26").Style = "Good"
$ws.Range("C26").Value = "ms"

# --- Row 27: 30 channels duration (formula now depends on RX_OFF/RX_ON) ---
$ws.Range("A27").Style = "Calculation"
$ws.Range("B27").Formula = "=(B25+B26)*30"
$ws.Range("B27").Style = "Calculation"

# --- Row 28: ON/OFF ratio ---
$ws.Range("A28").Value = "ON/OFF ratio"
$ws.Range("A28").Style = "Calculation"
$ws.Range("B28").Formula = "=100*B26/(B26+B25)"
$ws.Range("B28").Style = "Calculation"
$ws.Range("B28").NumberFormat = "0.0"
$ws.Range("C28").Value = "%"

# --- Row 22: Receive duration variant (*1.5) ---
$ws.Range("A22").Value = "Receive duration"
$ws.Range("B22").Formula = "=B21*1.5"
$ws.Range("B22").Style = "Normal"
$ws.Range("C22").Value = "ms"
$ws.Range("D22").Value = "*1.5"

# Widen column B slightly to fit the new numbers (~8.71 chars).
$ws.Columns.Item(2).ColumnWidth = 7.83

# Restore the selected cell to match the edited area.
$ws.Range("D23").Select()
